$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 4566.5835
$ws.Range("I19").Value = 4799.875
$ws.Range("J19").Value = 4100
$ws.Range("K19").Value = 4799.875
$ws.Range("L19").Value = 4100
$ws.Range("M19").Value = -4624.875
$ws.Range("N19").Value = -4450
$ws.Range("H32").Value = 4743.75
$ws.Range("I32").Value = 4500
$ws.Range("K32").Value = 4500
$ws.Range("M32").Value = -4174
$ws.Range("H64").Value = 8105.769
$ws.Range("I64").Value = 4098.75
$ws.Range("K64").Value = 4098.75
$ws.Range("M64").Value = -3850.75
$ws.Range("H67").Value = 8105.769
$ws.Range("I67").Value = 4098.75
$ws.Range("K67").Value = 4098.75
$ws.Range("M67").Value = -3240.75
$ws.Range("H75").Value = 19000
$ws.Range("J75").Value = 19000
$ws.Range("L75").Value = 19000
$ws.Range("N75").Value = -20872
$ws.Range("H76").Value = 6341.4165
$ws.Range("I76").Value = 4019.6
$ws.Range("J76").Value = 7999.857
$ws.Range("K76").Value = 4019.6
$ws.Range("L76").Value = 7999.857
$ws.Range("M76").Value = -3704.6
$ws.Range("N76").Value = -8629.857
$ws.Range("H78").Value = 19000
$ws.Range("J78").Value = 19000
$ws.Range("L78").Value = 57000
$ws.Range("N78").Value = -66360
$ws.Range("H79").Value = 6341.4165
$ws.Range("I79").Value = 4019.6
$ws.Range("J79").Value = 7999.857
$ws.Range("K79").Value = 4019.6
$ws.Range("L79").Value = 7999.857
$ws.Range("M79").Value = -2927.6
$ws.Range("N79").Value = -10183.857
$ws.Range("H106").Value = 2199
$ws.Range("I106").Value = 2447.8333
$ws.Range("K106").Value = 2447.8333
$ws.Range("M106").Value = -1816.8333
$ws.Range("H131").Value = 3364.6155
$ws.Range("J131").Value = 5851.25
$ws.Range("L131").Value = 17553.75
$ws.Range("N131").Value = -27633.75
$ws.Range("H132").Value = 5014.8184
$ws.Range("I132").Value = 1787.9231
$ws.Range("K132").Value = 5363.7693
$ws.Range("M132").Value = -2833.7693
$ws.Range("H137").Value = 32259836
$ws.Range("I137").Value = 50001172
$ws.Range("K137").Value = 150003516
$ws.Range("M137").Value = -150000966

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2007.4054
$ws.Range("I61").Value = 1182.5333
$ws.Range("K61").Value = 1182.5333
$ws.Range("M61").Value = -970.5333000000001
$ws.Range("H122").Value = 2007.3572
$ws.Range("I122").Value = 1563.909
$ws.Range("K122").Value = 4691.727000000001
$ws.Range("M122").Value = -2241.727000000001
$ws.Range("H132").Value = 4095.4
$ws.Range("I132").Value = 2816.3333
$ws.Range("K132").Value = 8448.999899999999
$ws.Range("M132").Value = -5918.999899999999
$ws.Range("H136").Value = 2007.4054
$ws.Range("I136").Value = 1182.5333
$ws.Range("K136").Value = 3547.5999
$ws.Range("M136").Value = -997.5999000000002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 60832.668
$ws.Range("J81").Value = 60832.668
$ws.Range("L81").Value = 60832.668
$ws.Range("N81").Value = -62954.668
$ws.Range("H84").Value = 60832.668
$ws.Range("J84").Value = 60832.668
$ws.Range("L84").Value = 182498.004
$ws.Range("N84").Value = -193106.004
$ws.Range("H107").Value = 7373
$ws.Range("I107").Value = 7322.2
$ws.Range("J107").Value = 7500
$ws.Range("K107").Value = 7322.2
$ws.Range("L107").Value = 7500
$ws.Range("M107").Value = -5402.2
$ws.Range("N107").Value = -11340

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3890.3809
$ws.Range("I31").Value = 1640.6522
$ws.Range("J31").Value = 6613.737
$ws.Range("K31").Value = 1640.6522
$ws.Range("L31").Value = 6613.737
$ws.Range("M31").Value = -1345.6522
$ws.Range("N31").Value = -7203.737
$ws.Range("H34").Value = 3890.3809
$ws.Range("I34").Value = 1640.6522
$ws.Range("J34").Value = 6613.737
$ws.Range("K34").Value = 1640.6522
$ws.Range("L34").Value = 6613.737
$ws.Range("M34").Value = -1438.6522
$ws.Range("N34").Value = -7017.737
$ws.Range("H58").Value = 2760
$ws.Range("I58").Value = 1428.9524
$ws.Range("J58").Value = 7418.6665
$ws.Range("K58").Value = 1428.9524
$ws.Range("L58").Value = 7418.6665
$ws.Range("M58").Value = -1225.9524
$ws.Range("N58").Value = -7824.6665
$ws.Range("H82").Value = 63374.715
$ws.Range("J82").Value = 63374.715
$ws.Range("L82").Value = 63374.715
$ws.Range("N82").Value = -64096.715
$ws.Range("H85").Value = 63374.715
$ws.Range("J85").Value = 63374.715
$ws.Range("L85").Value = 63374.715
$ws.Range("N85").Value = -65870.715
$ws.Range("H132").Value = 133338470
$ws.Range("I132").Value = 250003740
$ws.Range("K132").Value = 750011220
$ws.Range("M132").Value = -750008690
$ws.Range("H136").Value = 2760
$ws.Range("I136").Value = 1428.9524
$ws.Range("J136").Value = 7418.6665
$ws.Range("K136").Value = 4286.857199999999
$ws.Range("L136").Value = 22255.9995
$ws.Range("M136").Value = -1736.857199999999
$ws.Range("N136").Value = -27355.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 5112.0713
$ws.Range("I75").Value = 1501.4286
$ws.Range("K75").Value = 4504.2858
$ws.Range("M75").Value = -3506.2858
$ws.Range("H78").Value = 5112.0713
$ws.Range("I78").Value = 1501.4286
$ws.Range("K78").Value = 13512.8574
$ws.Range("M78").Value = -8520.857399999999
$ws.Range("H132").Value = 2312.8
$ws.Range("I132").Value = 2131.0908
$ws.Range("J132").Value = 2812.5
$ws.Range("K132").Value = 19179.8172
$ws.Range("L132").Value = 25312.5
$ws.Range("M132").Value = -16649.8172
$ws.Range("N132").Value = -30372.5
$ws.Range("H139").Value = 54741.26
$ws.Range("I139").Value = 57560.223
$ws.Range("K139").Value = 172680.669
$ws.Range("M139").Value = -167540.669

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 3048.5454
$ws.Range("I113").Value = 1183.5
$ws.Range("J113").Value = 4114.2856
$ws.Range("K113").Value = 1183.5
$ws.Range("L113").Value = 4114.2856
$ws.Range("M113").Value = 986.5
$ws.Range("N113").Value = -8454.285599999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 642.7692
$ws.Range("J22").Value = 644.75
$ws.Range("L22").Value = 644.75
$ws.Range("N22").Value = -1234.75
$ws.Range("H27").Value = 642.7692
$ws.Range("J27").Value = 644.75
$ws.Range("L27").Value = 644.75
$ws.Range("N27").Value = -858.75
$ws.Range("H40").Value = 5165.125
$ws.Range("J40").Value = 2000
$ws.Range("L40").Value = 2000
$ws.Range("N40").Value = -2272
$ws.Range("H122").Value = 4895
$ws.Range("I122").Value = 4333.6665
$ws.Range("K122").Value = 13000.9995
$ws.Range("M122").Value = -10550.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 125323280
$ws.Range("I4").Value = 840665
$ws.Range("J4").Value = 200012850
$ws.Range("K4").Value = 840665
$ws.Range("L4").Value = 200012850
$ws.Range("M4").Value = -840552
$ws.Range("N4").Value = -200013076
$ws.Range("H113").Value = 2979.9285
$ws.Range("I113").Value = 839.2
$ws.Range("K113").Value = 2517.6
$ws.Range("M113").Value = -347.6000000000004
$ws.Range("H132").Value = 3574985.5
$ws.Range("I132").Value = 4446941.5
$ws.Range("K132").Value = 13340824.5
$ws.Range("M132").Value = -13338294.5
